$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.082.93"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.216.64"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'241.31"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "'73.50"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.607"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "'43.12"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").Value = "'0.0954"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "'7.11"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "2.543.71"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "'14.21"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "'0.841"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "2.229.31"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "41.936.81"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  +9.51%  "
$ws.Range("D20").Value = "'72.74"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "'6.14"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'10.26"
$ws.Range("E22").Value = "  +17.94%  "
$ws.Range("D23").Value = "'229.69"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  -5.42%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'11.63"
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'3.60"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = "  +4.97%  "
$ws.Range("D30").Value = "'167.34"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").Value = "'20.60"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").Value = "'5.63"
$ws.Range("E32").Value = "  +8.22%  "
$ws.Range("D33").Value = "'0.0793"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'28.92"
$ws.Range("E35").Value = "  -5.08%  "
$ws.Range("E36").Value = "  -7.66%  "
$ws.Range("D37").Value = "'4.28"
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("D38").Value = "'0.0301"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("D39").Value = "'13.10"
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("B40").Value = "MultiversX"
$ws.Range("C40").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D40").Value = "'65.07"
$ws.Range("E40").Value = "  +4.06%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'2.12"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "'5.62"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").Value = "'0.199"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").Value = "'8.74"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").Value = "'104.51"
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'2.41"
$ws.Range("E47").Value = "  +6.30%  "
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'1.16"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "2.421.12"
$ws.Range("E51").Value = "  -1.38%  "
